# Simulated Wild Card round and logged it
# Appends the new playoff game's play-by-play yardage logs and bumps the
# aggregate season totals across the YDS / OFF / DEF / ST / TURNS / PEN sheets.

$wb = $excel.ActiveWorkbook

function Append-Tokens($ws, $cellRef, $tokens) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $cell.Value2 + " " + $tokens
}

function Add-Delta($ws, $cellRef, $delta) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $cell.Value2 + $delta
}

# ---------------------------------------------------------------------------
# YDS sheet: append this game's per-play yardage logs (Rushing/Passing,
# Offense/Defense).
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")
Append-Tokens $wsYDS "B2" "3 1 9 2 2 14 1 2 14 11 8 6 -1 2 -3 6 0 7 -1 4"
Append-Tokens $wsYDS "C2" "8 0 0 -3 14 -2 3 1 1 -2 0 23 1 2 7 5 3 15 2 1 0 11 7 5 8 -4 28 18 0 7 -1 7 10"
Append-Tokens $wsYDS "B3" "3 5 5 14 19 6 10 6 6 13 14 7 42 2 3 6 5 9 18 1 6 14 4 14 10 17 12 6"
Append-Tokens $wsYDS "C3" "8 44 0 4 12 12 18 4 4 30 6 2 -1 8 0 7 17 -5 11"

# ---------------------------------------------------------------------------
# OFF sheet: add this game's offensive down/distance + play counts.
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")
Add-Delta $wsOFF "C2" 14
Add-Delta $wsOFF "F2" 1
Add-Delta $wsOFF "G2" 4
Add-Delta $wsOFF "J2" 1
Add-Delta $wsOFF "N2" 3
Add-Delta $wsOFF "O2" 7
Add-Delta $wsOFF "P2" 6
Add-Delta $wsOFF "B3" 1
Add-Delta $wsOFF "C3" 19
Add-Delta $wsOFF "E3" 4
Add-Delta $wsOFF "F3" 16
Add-Delta $wsOFF "G3" 4
Add-Delta $wsOFF "H3" 2
Add-Delta $wsOFF "I3" 9
Add-Delta $wsOFF "J3" 6
Add-Delta $wsOFF "L3" 64
Add-Delta $wsOFF "M3" 34
Add-Delta $wsOFF "Q3" 91

# ---------------------------------------------------------------------------
# DEF sheet: add this game's defensive down/distance + play counts.
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")
Add-Delta $wsDEF "C2" 17
Add-Delta $wsDEF "D2" 3
Add-Delta $wsDEF "E2" 2
Add-Delta $wsDEF "F2" 3
Add-Delta $wsDEF "G2" 4
Add-Delta $wsDEF "H2" 1
Add-Delta $wsDEF "J2" 3
Add-Delta $wsDEF "N2" 3
Add-Delta $wsDEF "O2" 1
Add-Delta $wsDEF "P2" 1
Add-Delta $wsDEF "C3" 10
Add-Delta $wsDEF "E3" 1
Add-Delta $wsDEF "F3" 10
Add-Delta $wsDEF "G3" 4
Add-Delta $wsDEF "H3" 2
Add-Delta $wsDEF "I3" 7
Add-Delta $wsDEF "J3" 4
Add-Delta $wsDEF "L3" 36
Add-Delta $wsDEF "M3" 20
Add-Delta $wsDEF "Q3" 80

# ---------------------------------------------------------------------------
# ST sheet: bump special-teams counters and append this game's kickoff /
# punt distance logs (D / RA / RM per kick type).
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")
Add-Delta $wsST "B2" 6
Add-Delta $wsST "D2" 3
Add-Delta $wsST "F2" 3
Add-Delta $wsST "G2" 3
Add-Delta $wsST "H2" 1
Add-Delta $wsST "I2" 1
Add-Delta $wsST "L2" 1
Add-Delta $wsST "M2" 1
Add-Delta $wsST "N2" 1
Add-Delta $wsST "B3" 4

Append-Tokens $wsST "B4" "50 59"
Append-Tokens $wsST "B5" "18 22"
Append-Tokens $wsST "B6" "24 23"
Append-Tokens $wsST "D3" "55 51 56"
Append-Tokens $wsST "D4" "10 9 13"
Append-Tokens $wsST "D5" "4 0 0 11"

# ---------------------------------------------------------------------------
# TURNS sheet: add this game's road turnovers.
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")
Add-Delta $wsTURNS "B3" 1
Add-Delta $wsTURNS "E3" 2

# ---------------------------------------------------------------------------
# PEN sheet: add this game's penalty counts.
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")
Add-Delta $wsPEN "B2" 1
Add-Delta $wsPEN "D2" 1
Add-Delta $wsPEN "B3" 1
Add-Delta $wsPEN "D3" 3
Add-Delta $wsPEN "D4" 2
